$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.62'
$ws.Range("E2").Value = '''0.32%'
$ws.Range("E3").Value = '''2.82%'
$ws.Range("D4").Value = '''5.185'
$ws.Range("E4").Value = '''2.84%'
$ws.Range("D5").Value = '''0.05592'
$ws.Range("E5").Value = '''-0.05%'
$ws.Range("D6").Value = '''6.482'
$ws.Range("E6").Value = '''-1.28%'
$ws.Range("D7").Value = '''0.8132'
$ws.Range("E7").Value = '''-0.46%'
$ws.Range("E8").Value = '''1.31%'
$ws.Range("D9").Value = '''0.06921'
$ws.Range("E9").Value = '''-0.49%'
$ws.Range("D10").Value = '''0.02858'
$ws.Range("E10").Value = '''1.22%'
$ws.Range("D11").Value = '''0.09386'
$ws.Range("E11").Value = '''-0.11%'
$ws.Range("D12").Value = '''0.001512'
$ws.Range("E12").Value = '''-0.10%'
$ws.Range("B13").Value = 'TigerCash'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D13").Value = '''0.006136'
$ws.Range("E13").Value = '''0.25%'
$ws.Range("B14").Value = 'LEO'
$ws.Range("C14").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D14").Value = '''3.607'
$ws.Range("E14").Value = '''3.18%'
$ws.Range("B15").Value = 'GateToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D15").Value = '''3.028'
$ws.Range("E15").Value = '''0.50%'
$ws.Range("B16").Value = 'BTSEToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D16").Value = '''2.058'
$ws.Range("E16").Value = '''-1.58%'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D17").Value = '''0.009958'
$ws.Range("E17").Value = '''1,563.37%'
$ws.Range("D19").Value = '''0.1327'
$ws.Range("E19").Value = '''-0.72%'
$ws.Range("D20").Value = '''0.03122'
$ws.Range("E20").Value = '''-3.16%'
$ws.Range("D21").Value = '''0.1293'
$ws.Range("E21").Value = '''-1.96%'
$ws.Range("D22").Value = '''3.742'
$ws.Range("E22").Value = '''0.14%'
$ws.Range("D23").Value = '''0.04651'
$ws.Range("E23").Value = '''-0.94%'
$ws.Range("E25").Value = '''-0.10%'
$ws.Range("D26").Value = '''0.004548'
$ws.Range("E26").Value = '''5.95%'
$ws.Range("D27").Value = '''0.00009596'
$ws.Range("E27").Value = '''-1.01%'
$ws.Range("D28").Value = '''0.0001405'
$ws.Range("E28").Value = '''-27.54%'
$ws.Range("D40").Value = '''0.03648'
$ws.Range("E40").Value = '''-0.44%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006157'
$ws.Range("E41").Value = '''-0.62%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1053'
$ws.Range("E42").Value = '''-0.06%'
$ws.Range("D43").Value = '''0.002499'
$ws.Range("E43").Value = '''-3.83%'
$ws.Range("D44").Value = '''0.007982'
$ws.Range("E44").Value = '''7.59%'
$ws.Range("D45").Value = '''0.00005383'
$ws.Range("E45").Value = '''1.79%'
$ws.Range("D48").Value = '''0.002402'
$ws.Range("E48").Value = '''19.21%'
